# Auto-generated: apply scheduled market-data refresh to Chocobo Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 703.375
$ws.Range("I17").Value = 1100
$ws.Range("J17").Value = 604.21875
$ws.Range("K17").Value = 3300
$ws.Range("L17").Value = 1812.65625
$ws.Range("M17").Value = -3132
$ws.Range("N17").Value = -2148.65625
# Row 18
$ws.Range("H18").Value = 276.15384
$ws.Range("I18").Value = 217.27272
$ws.Range("J18").Value = 600
$ws.Range("K18").Value = 217.27272
$ws.Range("L18").Value = 600
$ws.Range("M18").Value = 66.72728000000001
$ws.Range("N18").Value = -1168
# Row 62
$ws.Range("H62").Value = 2822.1428
$ws.Range("J62").Value = 3868.5715
$ws.Range("L62").Value = 3868.5715
$ws.Range("N62").Value = -5116.5715
# Row 65
$ws.Range("H65").Value = 2822.1428
$ws.Range("J65").Value = 3868.5715
$ws.Range("L65").Value = 19342.8575
$ws.Range("N65").Value = -25582.8575
# Row 127
$ws.Range("H127").Value = 1512.4546
$ws.Range("I127").Value = 747.25
$ws.Range("J127").Value = 1949.7142
$ws.Range("K127").Value = 2241.75
$ws.Range("L127").Value = 5849.142599999999
$ws.Range("M127").Value = 2718.25
$ws.Range("N127").Value = -15769.1426
# Row 129
$ws.Range("H129").Value = 843.2033699999999
$ws.Range("J129").Value = 966.68085
$ws.Range("L129").Value = 2900.04255
$ws.Range("N129").Value = -12900.04255
# Row 132
$ws.Range("H132").Value = 18039316
$ws.Range("I132").Value = 20836762
$ws.Range("J132").Value = 1254633.4
$ws.Range("K132").Value = 62510286
$ws.Range("L132").Value = 3763900.2
$ws.Range("M132").Value = -62507756
$ws.Range("N132").Value = -3768960.2
# Row 135
$ws.Range("H135").Value = 502.1905
$ws.Range("I135").Value = 420.77777
$ws.Range("J135").Value = 990.6667
$ws.Range("K135").Value = 3786.99993
$ws.Range("L135").Value = 8916.0003
$ws.Range("M135").Value = -1251.99993
$ws.Range("N135").Value = -13986.0003
# Row 138
$ws.Range("H138").Value = 2091.5876
$ws.Range("I138").Value = 674.63416
$ws.Range("J138").Value = 3129
$ws.Range("K138").Value = 2023.90248
$ws.Range("L138").Value = 9387
$ws.Range("M138").Value = 3116.09752
$ws.Range("N138").Value = -19667

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2323.5
$ws.Range("I74").Value = 2156.6416
$ws.Range("J74").Value = 3586.8572
$ws.Range("K74").Value = 2156.6416
$ws.Range("L74").Value = 3586.8572
$ws.Range("M74").Value = -1282.6416
$ws.Range("N74").Value = -5334.8572
# Row 77
$ws.Range("H77").Value = 2323.5
$ws.Range("I77").Value = 2156.6416
$ws.Range("J77").Value = 3586.8572
$ws.Range("K77").Value = 10783.208
$ws.Range("L77").Value = 17934.286
$ws.Range("M77").Value = -6415.207999999999
$ws.Range("N77").Value = -26670.286
# Row 88
$ws.Range("H88").Value = 16668661
$ws.Range("I88").Value = 66666664
$ws.Range("J88").Value = 2660
$ws.Range("K88").Value = 66666664
$ws.Range("L88").Value = 2660
$ws.Range("M88").Value = -66666258
$ws.Range("N88").Value = -3472
# Row 91
$ws.Range("H91").Value = 16668661
$ws.Range("I91").Value = 66666664
$ws.Range("J91").Value = 2660
$ws.Range("K91").Value = 66666664
$ws.Range("L91").Value = 2660
$ws.Range("M91").Value = -66665260
$ws.Range("N91").Value = -5468
# Row 134
$ws.Range("H134").Value = 53332.668
$ws.Range("J134").Value = 53332.668
$ws.Range("L134").Value = 53332.668
$ws.Range("N134").Value = -63472.668

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1608.7646
$ws.Range("I86").Value = 1554.2142
$ws.Range("J86").Value = 1863.3334
$ws.Range("K86").Value = 1554.2142
$ws.Range("L86").Value = 1863.3334
$ws.Range("M86").Value = -431.2141999999999
$ws.Range("N86").Value = -4109.3334
# Row 89
$ws.Range("H89").Value = 1608.7646
$ws.Range("I89").Value = 1554.2142
$ws.Range("J89").Value = 1863.3334
$ws.Range("K89").Value = 7771.071
$ws.Range("L89").Value = 9316.666999999999
$ws.Range("M89").Value = -2155.071
$ws.Range("N89").Value = -20548.667
# Row 132
$ws.Range("H132").Value = 52375.24
$ws.Range("J132").Value = 52375.24
$ws.Range("L132").Value = 52375.24
$ws.Range("N132").Value = -62495.24

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 18520168
$ws.Range("I16").Value = 22223802
$ws.Range("J16").Value = 1993
$ws.Range("K16").Value = 22223802
$ws.Range("L16").Value = 1993
$ws.Range("M16").Value = -22223515
$ws.Range("N16").Value = -2567
# Row 113
$ws.Range("H113").Value = 18520168
$ws.Range("I113").Value = 22223802
$ws.Range("J113").Value = 1993
$ws.Range("K113").Value = 22223802
$ws.Range("L113").Value = 1993
$ws.Range("M113").Value = -22221632
$ws.Range("N113").Value = -6333
# Row 134
$ws.Range("H134").Value = 3159.2104
$ws.Range("I134").Value = 4357.607
$ws.Range("J134").Value = 2002.138
$ws.Range("K134").Value = 13072.821
$ws.Range("L134").Value = 6006.414
$ws.Range("M134").Value = -10537.821
$ws.Range("N134").Value = -11076.414

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1320.4884
$ws.Range("I5").Value = 325.21054
$ws.Range("J5").Value = 2108.4167
$ws.Range("K5").Value = 975.6316199999999
$ws.Range("L5").Value = 6325.250100000001
$ws.Range("M5").Value = -863.6316199999999
$ws.Range("N5").Value = -6549.250100000001
# Row 131
$ws.Range("H131").Value = 836.78
$ws.Range("I131").Value = 584.5
$ws.Range("J131").Value = 864.8111
$ws.Range("K131").Value = 1753.5
$ws.Range("L131").Value = 2594.4333
$ws.Range("M131").Value = 3286.5
$ws.Range("N131").Value = -12674.4333
# Row 133
$ws.Range("H133").Value = 3423.8235
$ws.Range("J133").Value = 2753.3333
$ws.Range("L133").Value = 8259.999899999999
$ws.Range("N133").Value = -18379.9999
# Row 135
$ws.Range("H135").Value = 1320.4884
$ws.Range("I135").Value = 325.21054
$ws.Range("J135").Value = 2108.4167
$ws.Range("K135").Value = 2926.89486
$ws.Range("L135").Value = 18975.7503
$ws.Range("M135").Value = -391.8948599999999
$ws.Range("N135").Value = -24045.7503
# Row 137
$ws.Range("H137").Value = 2474.2222
$ws.Range("I137").Value = 859.8182
$ws.Range("K137").Value = 2579.4546
$ws.Range("M137").Value = 2520.5454
# Row 140
$ws.Range("H140").Value = 3450.9
$ws.Range("I140").Value = 4802
$ws.Range("J140").Value = 2345.4546
$ws.Range("K140").Value = 14406
$ws.Range("L140").Value = 7036.3638
$ws.Range("M140").Value = -9226
$ws.Range("N140").Value = -17396.3638

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6426.054
$ws.Range("I70").Value = 5829.5415
$ws.Range("J70").Value = 7527.3076
$ws.Range("K70").Value = 5829.5415
$ws.Range("L70").Value = 7527.3076
$ws.Range("M70").Value = -5559.5415
$ws.Range("N70").Value = -8067.3076
# Row 73
$ws.Range("H73").Value = 6426.054
$ws.Range("I73").Value = 5829.5415
$ws.Range("J73").Value = 7527.3076
$ws.Range("K73").Value = 5829.5415
$ws.Range("L73").Value = 7527.3076
$ws.Range("M73").Value = -4893.5415
$ws.Range("N73").Value = -9399.3076
# Row 136
$ws.Range("H136").Value = 22227.666
$ws.Range("J136").Value = 22227.666
$ws.Range("L136").Value = 66682.99800000001
$ws.Range("N136").Value = -71782.99800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 30389
$ws.Range("J92").Value = 30389
$ws.Range("L92").Value = 30389
$ws.Range("N92").Value = -35381

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 10754427
$ws.Range("I132").Value = 684.7826
$ws.Range("K132").Value = 2054.3478
$ws.Range("M132").Value = 475.6522
